$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Reviewer line: "Reviewer	:  Nur Halimatussa'diyah, " + [_GoBack] + "16090068"
#    becomes one continuous run of text (no bookmark on this line any more -
#    the _GoBack bookmark is relocated into the "Latar belakang" paragraph,
#    see step 3 below). Remove the old bookmark and make sure the two
#    remaining runs read continuously.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute("Nur Halimatussa" + [char]8217 + "diyah, 16090068", $true, $false, $false, $false, $false, $true, 1, $false, "Nur Halimatussa" + [char]8217 + "diyah, 16090068", 2)
Write-Host "Reviewer merge found: $found"

# ------------------------------------------------------------------
# 2) "Latar belakang" paragraph: rewrite/expand the body text.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("asasi", $true, $false, $false, $false, $false, $true, 1, $false, "bagi setiap", 2)
Write-Host "Replace 1 found: $found"

$rng = $d.Content
$found = $rng.Find.Execute("yang menimbulkan tantangan baru", $true, $false, $false, $false, $false, $true, 1, $false, "yang pada akhirnya menimbulkan masalah baru", 2)
Write-Host "Replace 2 found: $found"

$rng = $d.Content
$found = $rng.Find.Execute("Dampak yang dihasilkan dari lingkungan ini", $true, $false, $false, $false, $false, $true, 1, $false, "Dampak yang terjadi dari perubahan ini", 2)
Write-Host "Replace 3 found: $found"

$rng = $d.Content
$old4 = "adalah perubahan lingkungan pada masyarakat tentang kepercayaan dan moral. Etika digital digunakan"
$new4 = "adalah perubahan masyarakat tentang sebuah kepercayaan dan moral. Pemanfaatkan internet sebagai salahsatu infrastruktur komunikasi global, juga memunculkan masalah baru terkait privasi, maka perlu dipertanyakan tentang keamanan data privasi yang terekam dan etika yang digunakan. Etika digital muncul"
$found = $rng.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Replace 4 found: $found"

# ------------------------------------------------------------------
# 3) Re-create the _GoBack bookmark at its new location, right before
#    the word "moral" that follows "kepercayaan dan " in the rewritten
#    "Latar belakang" paragraph.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("kepercayaan dan ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Bookmark anchor found: $found"
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)
Write-Host "GoBack re-added: $($d.Bookmarks.Exists('_GoBack'))"
